{"js": "// Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line,\n// the copyright/footer line that follows it, and the blank paragraph that\n// separates them from the \"Requisitos\" section above, while leaving the\n// blank paragraph (and page-break paragraph) that originally trailed the\n// footer text untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\n// Locate the \"Requisitos\" entry paragraph that precedes the block to remove.\nconst anchorMarker = \"LOQ4044\";\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorMarker) !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph containing '\" + anchorMarker + \"'\");\n}\n\n// The three paragraphs immediately after the anchor are:\n//   1) an empty spacer paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) the \"\u00a9 2020 ... Creative Commons Attribution\" footer line\n// Delete them in reverse order so earlier deletions don't shift later indexes.\nfor (let offset = 3; offset >= 1; offset--) {\n  paragraphs.items[anchorIndex + offset].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line,\n# the copyright/footer line that follows it, and the blank paragraph that\n# separates them from the \"Requisitos\" section above, while leaving the\n# blank paragraph (and page-break paragraph) that originally trailed the\n# footer text untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Requisitos\" entry paragraph that precedes the block to remove.\n$anchorMarker = \"LOQ4044\"\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($text -like \"*$anchorMarker*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph containing '$anchorMarker'\"\n}\n\n# The three paragraphs immediately after the anchor are:\n#   1) an empty spacer paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) the \"(c) 2020 ... Creative Commons Attribution\" footer line\n# Delete them starting from the furthest one so earlier deletions don't shift\n# the indices of paragraphs not yet removed.\nfor ($offset = 3; $offset -ge 1; $offset--) {\n    $d.Paragraphs.Item($anchorIndex + $offset).Range.Delete()\n}\n"}
